{"js": "// Better handle copying paragraph styles:\n// Whenever a paragraph uses the \"MSC_Join\" (styleId \"MSCJoin\") style, the\n// paragraph immediately following it (normally holding the \"[...]\" elision\n// marker) and the empty paragraph after that should inherit the same\n// \"MSC_Join\" paragraph style, instead of being left on the default style.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\nconst JOIN_STYLE = \"MSC_Join\";\nconst items = paragraphs.items;\n\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style !== JOIN_STYLE) continue;\n\n  // The \"[...]\" paragraph right after the MSC_Join paragraph.\n  const ellipsisPara = items[i + 1];\n  if (!ellipsisPara || ellipsisPara.text !== \"[...]\") continue;\n\n  // The empty paragraph right after the \"[...]\" paragraph.\n  const blankPara = items[i + 2];\n\n  if (ellipsisPara.style !== JOIN_STYLE) {\n    ellipsisPara.style = JOIN_STYLE;\n  }\n  if (blankPara && blankPara.text === \"\" && blankPara.style !== JOIN_STYLE) {\n    blankPara.style = JOIN_STYLE;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Better handle copying paragraph styles:\n# Whenever a paragraph uses the \"MSC_Join\" (styleId \"MSCJoin\") style, the\n# paragraph immediately following it (normally holding the \"[...]\" elision\n# marker) and the empty paragraph after that should inherit the same\n# \"MSC_Join\" paragraph style, instead of being left on the default style.\n\n$d = $word.ActiveDocument\n$joinStyle = \"MSC_Join\"\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -ne $joinStyle) {\n        continue\n    }\n\n    # The \"[...]\" paragraph right after the MSC_Join paragraph.\n    if ($i + 1 -gt $count) {\n        continue\n    }\n    $ellipsisPara = $d.Paragraphs.Item($i + 1)\n    if ($ellipsisPara.Range.Text.Trim() -ne \"[...]\") {\n        continue\n    }\n\n    if ($ellipsisPara.Style.NameLocal -ne $joinStyle) {\n        $ellipsisPara.Style = $joinStyle\n    }\n\n    # The empty paragraph right after the \"[...]\" paragraph.\n    if ($i + 2 -le $count) {\n        $blankPara = $d.Paragraphs.Item($i + 2)\n        if ($blankPara.Style.NameLocal -ne $joinStyle) {\n            $blankPara.Style = $joinStyle\n        }\n    }\n}\n"}
